$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.879.77'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.265.73'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.534'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.484'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '54.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '32.28'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0800'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.87%  '
$ws.Range('E13').Value = '  +2.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.616.80'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('E16').Value = '  +2.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.266.97'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.757'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.786.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0903'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '241.66'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.26%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  +4.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.94'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.14%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.27%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.08'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.85'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.51%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.16'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0746'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.77%  '
$ws.Range('E36').Value = '  +3.01%  '
$ws.Range('E37').Value = '  +2.45%  '
$ws.Range('E38').Value = '  +5.17%  '
$ws.Range('E39').Value = '  +3.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.64%  '
$ws.Range('E41').Value = '  +4.55%  '
$ws.Range('E42').Value = '  +6.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.071.51'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.46'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.29%  '
$ws.Range('E45').Value = '  +2.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.10'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.17%  '
$ws.Range('E48').Value = '  +2.02%  '
$ws.Range('E49').Value = '  +3.96%  '
$ws.Range('E50').Value = '  +3.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.73'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.67%  '
